$d = $word.ActiveDocument
$dash = [char]0x2013

# 1. GPA: 3.84/4.00 -> 3.84
$p6 = $d.Paragraphs.Item(6).Range
$p6.Find.Execute(".84/4.00", $false, $false, $false, $false, $false, $true, 1, $false, ".84", 2)

# 2. "June 2021 - Present" -> "June 2021 <en-dash> August 2021" (Software Developer Quality Management Intern role)
$p9 = $d.Paragraphs.Item(9).Range
$p9.Find.Execute("June 2021 - Present", $false, $false, $false, $false, $false, $true, 1, $false, "June 2021 $dash August 2021", 2)

# 3. Rewrite the "Developing..." bullet
$p10 = $d.Paragraphs.Item(10).Range
$p10.Find.Execute("Developing the CI/CD dashboard for enterprise products", $false, $false, $false, $false, $false, $true, 1, $false, "Created the CI/CD dashboard for the Spectrum Spatial enterprise product", 2)
$p10b = $d.Paragraphs.Item(10).Range
$p10b.Find.Execute("using Python by converting non-useable test data into readable display formats for the dashboard rendering system.", $false, $false, $false, $false, $false, $true, 1, $false, "by using a Python Gitlab workflow to convert Junit, Nunit, and TestNG tests into useable Elasticsearch documents.", 2)

# 4. The "Adding 3 new test formats..." bullet is truncated to "Re"
$p11 = $d.Paragraphs.Item(11).Range
$rng11 = $d.Range($p11.Start, $p11.End)
$rng11.Text = "Re"

# 5/6. "May 2021 - Present" -> "May 2021 <en-dash> August 2021" (Software Developer, University of Toronto role)
$p13 = $d.Paragraphs.Item(13).Range
$p13.Find.Execute("May 2021 - Present", $false, $false, $false, $false, $false, $true, 1, $false, "May 2021 $dash August 2021", 2)

# 7. "Collaborating on the development team for " -> "Worked on the development team for "
$p14 = $d.Paragraphs.Item(14).Range
$p14.Find.Execute("Collaborating on the development team for", $false, $false, $false, $false, $false, $true, 1, $false, "Worked on the development team for", 2)

# 8. "Updating documentation..." -> "Updated documentation..."
$p15 = $d.Paragraphs.Item(15).Range
$p15.Find.Execute("Updating documentation", $false, $false, $false, $false, $false, $true, 1, $false, "Updated documentation", 2)

# 9. "Updating error report page" -> "Updated the error report page"
$p16 = $d.Paragraphs.Item(16).Range
$p16.Find.Execute("Updating error report page", $false, $false, $false, $false, $false, $true, 1, $false, "Updated the error report page", 2)
